$wb = $excel.ActiveWorkbook

# Grab the two existing sheets by their current names.
$sheetA = $wb.Worksheets.Item("hotel_info")   # physically sheet1.xml (rId1)
$sheetB = $wb.Worksheets.Item("review_info")  # physically sheet2.xml (rId2)

# Swap their names (hotel_info <-> review_info) via a temporary name to
# avoid a naming collision. $sheetA keeps pointing at the same physical
# worksheet (sheet1.xml) and $sheetB at sheet2.xml throughout.
$sheetA.Name = "tmp_swap_name"
$sheetB.Name = "hotel_info"
$sheetA.Name = "review_info"

# --- $sheetA is now named "review_info" (sheet1.xml / rId1) ---------------
# Replace its single header row with the full review_info header list and
# make sure there is no leftover data row.
$reviewHeaders = @("STR","reviewer_ID","reviewer_name","Review_ID","Date_of_scraping","ReviewURL","Tripadvisor_gcode","Tripadvisor_dcode","Tripadvisor_rcode","review_date","review_title","review_content","review_rating","trip_month","trip_purpose","value","rooms","Location","Cleanliness","Sleep Quality","Service","Picture(yes=1)","respondent","response_date","response_text")

$sheetA.Cells.Clear()
for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $sheetA.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- $sheetB is now named "hotel_info" (sheet2.xml / rId2) ----------------
# Replace its header row with the hotel_info headers (now including the
# new "State" column) and write the corresponding data row.
$hotelHeaders = @("STR","Hotel_Name","State","City","Zip","TA_ReviewURL","Tripadvisor_Hotel_Name","English_Reviews_num","Local_Rank","Total_Reviews_num")

$sheetB.Cells.Clear()
for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $sheetB.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

$sheetB.Cells.Item(2, 1).Value = 37670
$sheetB.Cells.Item(2, 2).Value = "Holiday Inn Express New Orleans Downtown"
$sheetB.Cells.Item(2, 3).Value = "Louisiana"
$sheetB.Cells.Item(2, 4).Value = "New Orleans"
$sheetB.Cells.Item(2, 5).Value = 70112
$sheetB.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d93147-Reviews-Holiday_Inn_Express_New_Orleans_Downtown-New_Orleans_Louisiana.html"
$sheetB.Cells.Item(2, 7).Value = "Holiday Inn Express New Orleans Downtown"
# These three are digit strings that must stay text (shared-string), not
# numbers, so force them with a leading apostrophe, then strip the
# resulting "quote prefix" formatting so the cell keeps the default style.
$sheetB.Cells.Item(2, 8).Value = "'259"
$sheetB.Cells.Item(2, 8).ClearFormats()
$sheetB.Cells.Item(2, 9).Value = "'124"
$sheetB.Cells.Item(2, 9).ClearFormats()
$sheetB.Cells.Item(2, 10).Value = "'288"
$sheetB.Cells.Item(2, 10).ClearFormats()
